# Added a new Test case
#
# Sheet1 holds a small table of test-case rows (columns A:M, header in row 1).
# This change:
#   1. Renames the existing row-4 test case from "Framework_002" to
#      "Framework_003" (it was a duplicate of row 3's name).
#   2. Adds a brand new row 5 ("Framework_004") that otherwise mirrors row 4,
#      but with a different person (Pratik Sharma / new address).
#   3. Extends the two list data-validations down to the new row.
#   4. Adds the matching mailto hyperlinks for the new row's C and M cells.
#   5. Leaves the grid with N9 selected, as last left by the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the existing row 4 test case name -------------------------
$ws.Range("A4").Value = "Framework_003"

# --- 2. Hyperlink row 5's future Password/Email cells while they are still
#        blank, then duplicate row 4 over the top. Copying row 4 afterwards
#        re-applies row 4's exact cell styles (incl. the Hyperlink style) to
#        row 5, so the new row ends up visually identical to row 4.
[void]$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Test@123")
[void]$ws.Hyperlinks.Add($ws.Range("M5"), "mailto:tooolsqa@gmail.com")

$ws.Range("A4:M4").Copy($ws.Range("A5:M5"))

# --- 3. Edit the handful of cells that differ for the new test case ------
$ws.Range("A5").Value = "Framework_004"
$ws.Range("G5").Value = "Pratik"
$ws.Range("I5").Value = "29, LimeSquare, City Road"

# --- 4. Extend the list validations on E and F down through row 5 --------
$ws.Range("E2:E4").Validation.Delete()
$ws.Range("F2:F4").Validation.Delete()
$ws.Range("E2:E5").Validation.Add(3, 1, 1, '"Accessories, iMacs, iPads, iPhones"')
$ws.Range("F2:F5").Validation.Add(3, 1, 1, '"Product 1, Product 2, Product 3, Product 4"')

# --- 5. Leave the same cell selected as in the authored workbook ---------
[void]$ws.Range("N9").Select()
